# The deck ships two OOXML theme parts:
#   theme1.xml -> "Office Theme" colour scheme (only ever linked from the
#                 Notes Master, which PowerPoint's object model does not
#                 expose as a separately themeable master here)
#   theme2.xml -> "Integral" colour scheme (linked from the Slide Master,
#                 i.e. the theme that actually drives the look of every
#                 slide/layout in the deck)
#
# The commit swaps the two themes' colour schemes. The part of that swap
# that is visible/reachable through the PowerPoint object model is the
# Slide Master's theme, so re-point every slot of its
# Theme.ThemeColorScheme from the old "Integral" palette to the
# "Office Theme" palette (the colours theme1.xml carried before the edit).
#
# ThemeColorScheme indices (1-based), matching OOXML <a:clrScheme> order:
#   1 dk1  2 lt1  3 dk2  4 lt2  5 accent1  6 accent2  7 accent3
#   8 accent4  9 accent5  10 accent6  11 hlink  12 folHlink

$p = $ppt.ActivePresentation
$cs = $p.SlideMaster.Theme.ThemeColorScheme

$cs.Item(1).RGB  = 0          # dk1      000000
$cs.Item(2).RGB  = 16777215   # lt1      FFFFFF
$cs.Item(3).RGB  = 6968388    # dk2      44546A
$cs.Item(4).RGB  = 15132391   # lt2      E7E6E6
$cs.Item(5).RGB  = 13998939   # accent1  5B9BD5
$cs.Item(6).RGB  = 3243501    # accent2  ED7D31
$cs.Item(7).RGB  = 10855845   # accent3  A5A5A5
$cs.Item(8).RGB  = 49407      # accent4  FFC000
$cs.Item(9).RGB  = 12874308   # accent5  4472C4
$cs.Item(10).RGB = 4697456    # accent6  70AD47
$cs.Item(11).RGB = 12673797   # hlink    0563C1
$cs.Item(12).RGB = 7491477    # folHlink 954F72
